# Weekly update: insert 3 new rows of fresh Espárragos price data at the
# top of the date-ordered block (row 27), shifting the existing rows
# (27-59) down to (30-62).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new blank rows starting at row 27 (old rows 27-59 shift to 30-62).
$ws.Range("A27:A29").EntireRow.Insert()

# Data shared by every row in this subset (constant across the whole sheet).
$mercadoId = 9
$mercado = "Vega Central Mapocho de Santiago"
$region = "Metropolitana"
$codreg = 13
$categoriaId = 300000000
$categoria = "Espárragos"
$clasificacion = "Hortaliza"

# New row 27: Banquete
$ws.Cells.Item(27, 1).Value = $mercadoId
$ws.Cells.Item(27, 2).Value = $mercado
$ws.Cells.Item(27, 3).Value = $region
$ws.Cells.Item(27, 4).Value = "2021-10-14"
$ws.Cells.Item(27, 5).Value = $codreg
$ws.Cells.Item(27, 6).Value = $categoriaId
$ws.Cells.Item(27, 7).Value = $categoria
$ws.Cells.Item(27, 8).Value = "Sin especificar"
$ws.Cells.Item(27, 9).Value = "Banquete"
$ws.Cells.Item(27, 10).Value = 160
$ws.Cells.Item(27, 11).Value = 1500
$ws.Cells.Item(27, 12).Value = 1500
$ws.Cells.Item(27, 13).Value = 1500
$ws.Cells.Item(27, 14).Value = "$/kilo"
$ws.Cells.Item(27, 15).Value = "Provincia de Linares"
$ws.Cells.Item(27, 16).Value = 1500
$ws.Cells.Item(27, 17).Value = 1
$ws.Cells.Item(27, 18).Value = $clasificacion

# New row 28: Primera
$ws.Cells.Item(28, 1).Value = $mercadoId
$ws.Cells.Item(28, 2).Value = $mercado
$ws.Cells.Item(28, 3).Value = $region
$ws.Cells.Item(28, 4).Value = "2021-10-14"
$ws.Cells.Item(28, 5).Value = $codreg
$ws.Cells.Item(28, 6).Value = $categoriaId
$ws.Cells.Item(28, 7).Value = $categoria
$ws.Cells.Item(28, 8).Value = "Sin especificar"
$ws.Cells.Item(28, 9).Value = "Primera"
$ws.Cells.Item(28, 10).Value = 430
$ws.Cells.Item(28, 11).Value = 1300
$ws.Cells.Item(28, 12).Value = 1300
$ws.Cells.Item(28, 13).Value = 1300
$ws.Cells.Item(28, 14).Value = "$/kilo"
$ws.Cells.Item(28, 15).Value = "Provincia de Linares"
$ws.Cells.Item(28, 16).Value = 1300
$ws.Cells.Item(28, 17).Value = 1
$ws.Cells.Item(28, 18).Value = $clasificacion

# New row 29: Segunda
$ws.Cells.Item(29, 1).Value = $mercadoId
$ws.Cells.Item(29, 2).Value = $mercado
$ws.Cells.Item(29, 3).Value = $region
$ws.Cells.Item(29, 4).Value = "2021-10-14"
$ws.Cells.Item(29, 5).Value = $codreg
$ws.Cells.Item(29, 6).Value = $categoriaId
$ws.Cells.Item(29, 7).Value = $categoria
$ws.Cells.Item(29, 8).Value = "Sin especificar"
$ws.Cells.Item(29, 9).Value = "Segunda"
$ws.Cells.Item(29, 10).Value = 250
$ws.Cells.Item(29, 11).Value = 1100
$ws.Cells.Item(29, 12).Value = 1100
$ws.Cells.Item(29, 13).Value = 1100
$ws.Cells.Item(29, 14).Value = "$/kilo"
$ws.Cells.Item(29, 15).Value = "Provincia de Linares"
$ws.Cells.Item(29, 16).Value = 1100
$ws.Cells.Item(29, 17).Value = 1
$ws.Cells.Item(29, 18).Value = $clasificacion
